$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new result row (row 38) produced by the latest script run.
# The date must stay as plain text (matches existing rows, all stored as
# text dates), so we use a leading apostrophe to stop Excel from turning
# it into a serial date number, then reset the style back to the sheet's
# default "Normal" so no stray number-format / quote-prefix style sticks.
$ws.Range("A38").Value = "'2025-04-01"
$ws.Range("A38").Style = "Normal"

$ws.Range("B38").Value = "zone tampon"
$ws.Range("C38").Value = 50
$ws.Range("D38").Value = 1
